# ---------------------------------------------------------------------------
# Edit script: splits two paragraphs' runs to match the target OOXML diff.
#
# Change 1: "  My public key is 17, 23." run gets split into
#           "  My public key is 17, " | "1" | "23."
#           (net effect: the visible text becomes "...17, 123.")
#
# Change 2: the base64 blobs in the "Hint 1" / "Hint 2" paragraphs (at the
#           very end of the document) are swapped, and each paragraph's
#           text gets split across multiple runs.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

function Split-AtBoldToggle($range) {
    # Toggling Bold on/off (ending back at its original value) forces the
    # engine to re-tokenize runs at this range's boundaries without
    # altering the effective formatting.
    $range.Bold = 1
    $range.Bold = 0
}

# ===========================================================================
# Change 1: "My public key is 17, 23."
# ===========================================================================

$rng = $d.Content
$found = $rng.Find.Execute("17, 23.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find '17, 23.' text"
}
$matchStart = $rng.Start
$matchEnd = $rng.End

# Insert the extra "1" right before "23." (after "17, ").
$insertPoint = $d.Range($matchStart + 4, $matchStart + 4)
$insertPoint.InsertBefore("1")
$newMatchEnd = $matchEnd + 1

# Split boundaries:
#   [matchStart,      matchStart + 5) -> "...17, 1"  (kept together with the
#                                          preceding "  My public key is " text)
#   [matchStart + 5,  newMatchEnd)    -> "23."
$splitAfter1 = $matchStart + 5

$runWithInsertedOne = $d.Range($matchStart, $splitAfter1)
Split-AtBoldToggle $runWithInsertedOne

$runTrailing = $d.Range($splitAfter1, $newMatchEnd)
Split-AtBoldToggle $runTrailing

# That produced a split after the inserted "1" character; we additionally
# need a split right after "17, " (before the inserted "1") so that "1"
# becomes its own run.
$splitBefore1 = $matchStart + 4
$runJustOne = $d.Range($splitBefore1, $splitAfter1)
Split-AtBoldToggle $runJustOne

# ===========================================================================
# Change 2: Hint 1 / Hint 2 paragraphs near the end of the document
# ===========================================================================

$cybermenB64 = "VGhlIEN5YmVybWVuIGhhdmUgY2hvc2VuIGEgbW9kdWx1cyAoMjExKSBhbmQgY3VydmUgdGhhdCByZXN1bHRzIGluIDIzMiBwb2ludHMuICBJIGd1ZXNzIHlvdSBjb3VsZCBnZXQgYWxsIDIzMiBwb2ludHMgZnJvbSB0aGUgY2FsY3VsYXRvciBhbmQgdGhlbiBicmVhayBvbmUgb2YgdGhlIHB1YmxpYyBrZXlzIHRvIGdldCBhIHByaXZhdGUga2V5LiAgVGhhdCB3b3VsZCBiZSBhIHJlYWwgcGFpbiwgdGhvdWdoLiAgVGhlcmUgaGFzIGdvdCB0byBiZSBhbiBlYXNpZXIgd2F5Lg=="
$beforeYouTryB64 = "QmVmb3JlIHlvdSB0cnkgdG8gYnJlYWsgYW55b25l4oCZcyBrZXksIHB1dCB0aGVpciBjdXJ2ZSwgZmllbGQsIGFuZCBiYXNlIHBvaW50IFAgaW50byB0aGUgY2FsY3VsYXRvciBhbmQgcGxheSB3aXRoIE4uICBZb3UgbWF5IGZpbmQgdGhlcmUgYXJlIG5vdCBtYW55IGNob2ljZXMgZm9yIHRoZSBzaGFyZWQga2V5Lg=="
$cybermenB64minus1 = $cybermenB64.Substring(0, $cybermenB64.Length - 1)

# Locate the two target paragraphs robustly by content, since the doc has
# other "Hint 1" / "Hint 2" paragraphs elsewhere.
$cnt = $d.Paragraphs.Count
$hint1Idx = -1
$hint2Idx = -1
for ($i = 1; $i -le $cnt; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text
    if ($t.StartsWith("Hint 1 VGhl")) {
        $hint1Idx = $i
    }
    if ($hint1Idx -ne -1 -and $i -eq ($hint1Idx + 1) -and $t.StartsWith("Hint 2")) {
        $hint2Idx = $i
    }
}
if ($hint1Idx -eq -1 -or $hint2Idx -eq -1) {
    throw "Could not locate Hint 1/Hint 2 paragraphs"
}

# --- Swap the base64 payloads between the two paragraphs ---

$p1 = $d.Paragraphs($hint1Idx)
$r1 = $p1.Range
$r1.Find.Execute($cybermenB64, $true, $false, $false, $false, $false, $true, 1, $false, ($beforeYouTryB64 + " ="), 2)

$p2 = $d.Paragraphs($hint2Idx)
$r2 = $p2.Range
$r2.Find.Execute($beforeYouTryB64, $true, $false, $false, $false, $false, $true, 1, $false, $cybermenB64minus1, 2)

# --- Re-split the "Hint 1" paragraph into 4 runs ---
# "Hint 1 " | beforeYouTryB64 | " " | "="

$p1 = $d.Paragraphs($hint1Idx)
$p1Start = $p1.Range.Start

$lenPrefix = "Hint 1 ".Length
$lenB64 = $beforeYouTryB64.Length
$split1 = $p1Start + $lenPrefix
$split2 = $split1 + $lenB64
$split3 = $split2 + 1

Split-AtBoldToggle ($d.Range($p1Start, $split1))
Split-AtBoldToggle ($d.Range($split1, $split2))
Split-AtBoldToggle ($d.Range($split2, $split3))

# --- Re-split the "Hint 2" paragraph so <w:br/> ends its own run ---
# "Hint 2" | <w:br/> | cybermenB64minus1

$p2 = $d.Paragraphs($hint2Idx)
$p2Start = $p2.Range.Start
$brPos = $p2Start + "Hint 2".Length
$afterBrPos = $brPos + 1

Split-AtBoldToggle ($d.Range($p2Start, $afterBrPos))

Write-Host "Edits applied."
